$d = $word.ActiveDocument

# Locate the paragraph that ends with "...class to make text white in Bootstrap"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "class to make text white in Bootstrap") {
        $target = $p
    }
}

# Insert a new paragraph right after it, then set its text/content.
$target.Range.InsertParagraphAfter()
$newPara = $target.Next()
$newPara.Range.Text = "In react we use {} when we want to write js in between the tags."
